$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the split runs in the "Value for Matrix A is (matrix name, ..."
#    sentence into a single run with the same combined text.
# ---------------------------------------------------------------------
$mergeText = " is (matrix name, column number, value), for value 1, the value would be (A, 1, 1)."
$rng = $d.Content
$rng.Find.Execute($mergeText, $true, $false, $false, $false, $false, $true, 1, $false, $mergeText, 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Remove the stray _GoBack bookmark that sits by itself in the blank
#    paragraph above the "Characteristics of Iliad and Odyssey" heading.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# 3. Bold the "In Odyssey" lead-in of the Odyssey bullet.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("In Odyssey, there are words like home", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$odysseyBold = $d.Range($rng.Start, $rng.Start + 10)
$odysseyBold.Font.Bold = 1

# ---------------------------------------------------------------------
# 4. Bold the "In Iliad" lead-in of the Iliad bullet and wrap it with a
#    fresh _GoBack bookmark (this is where the user's cursor last was).
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("In Iliad, there are words like war", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$iliadBold = $d.Range($rng.Start, $rng.Start + 8)
$iliadBold.Font.Bold = 1
$d.Bookmarks.Add("_GoBack", $iliadBold)

# ---------------------------------------------------------------------
# 5. Add a new bullet after the Iliad characteristics bullet.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("seems like it is related to war or anger.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newPara.Range.InsertBefore("There are not many things common in the text.")
